# Force Sensor Testing.xlsx - add two new sheets with FSR + sensor retest data
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Tweak Sheet1's current selection (cosmetic, matches author's last click)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2:B2").Select()

# ------------------------------------------------------------------
# 2. Create the two new worksheets at the end of the workbook, in order
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "FSR_2 Input Test + Sensor 2"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "FSR_1 Input Test + Sensor 1"

# ------------------------------------------------------------------
# 3. Populate "FSR_2 Input Test + Sensor 2" (sheet3)
# ------------------------------------------------------------------
$ws3.Range("A1").Value = "FSR_2 Input Test + Sensor 2"

# Header row - written in the exact order the author originally typed them
$ws3.Range("A2").Value = "Load (kg)"
$ws3.Range("C2").Value = "ADC_test 1"
$ws3.Range("B2").Value = "Voltage_test 1 (V)"
$ws3.Range("D2").Value = "Voltage_test 2 (V)"
$ws3.Range("E2").Value = "ADC_test 2"
$ws3.Range("F2").Value = "Voltage_test 3 (V)"
$ws3.Range("G2").Value = "ADC_test 3"
$ws3.Range("H2").Value = "Voltage_test 4 (V)"
$ws3.Range("I2").Value = "ADC_test 4"

# Data rows 3-12
$rows = New-Object 'object[,]' 10,9
$rows[0,0]=0.1;  $rows[0,1]=0.14;  $rows[0,2]=48;  $rows[0,3]=0.12;  $rows[0,4]=45;  $rows[0,5]=0.13;  $rows[0,6]=30;  $rows[0,7]=0.129; $rows[0,8]=25
$rows[1,0]=0.2;  $rows[1,1]=0.19;  $rows[1,2]=80;  $rows[1,3]=0.15;  $rows[1,4]=46;  $rows[1,5]=0.19;  $rows[1,6]=109; $rows[1,7]=0.145; $rows[1,8]=45
$rows[2,0]=0.3;  $rows[2,1]=0.27;  $rows[2,2]=175; $rows[2,3]=0.24;  $rows[2,4]=145; $rows[2,5]=0.19;  $rows[2,6]=100; $rows[2,7]=0.21;  $rows[2,8]=140
$rows[3,0]=0.4;  $rows[3,1]=0.33;  $rows[3,2]=260; $rows[3,3]=0.28;  $rows[3,4]=210; $rows[3,5]=0.21;  $rows[3,6]=120; $rows[3,7]=0.23;  $rows[3,8]=150
$rows[4,0]=0.5;  $rows[4,1]=0.39;  $rows[4,2]=330; $rows[4,3]=0.26;  $rows[4,4]=290; $rows[4,5]=0.25;  $rows[4,6]=150; $rows[4,7]=0.26;  $rows[4,8]=190
$rows[5,0]=1;    $rows[5,1]=0.58;  $rows[5,2]=560; $rows[5,3]=0.5;   $rows[5,4]=450; $rows[5,5]=0.32;  $rows[5,6]=320; $rows[5,7]=0.4;   $rows[5,8]=360
$rows[6,0]=1.5;  $rows[6,1]=0.56;  $rows[6,2]=549; $rows[6,3]=0.55;  $rows[6,4]=540; $rows[6,5]=0.55;  $rows[6,6]=530; $rows[6,7]=0.5;   $rows[6,8]=460
$rows[7,0]=2;    $rows[7,1]=0.6;   $rows[7,2]=590; $rows[7,3]=0.64;  $rows[7,4]=650; $rows[7,5]=0.61;  $rows[7,6]=590; $rows[7,7]=0.56;  $rows[7,8]=560
$rows[8,0]=2.5;  $rows[8,1]=0.765; $rows[8,2]=790; $rows[8,3]=0.68;  $rows[8,4]=680; $rows[8,5]=0.71;  $rows[8,6]=720; $rows[8,7]=0.63;  $rows[8,8]=650
$rows[9,0]=3;    $rows[9,1]=0.82;  $rows[9,2]=860; $rows[9,3]=0.71;  $rows[9,4]=720; $rows[9,5]=0.76;  $rows[9,6]=790; $rows[9,7]=0.71;  $rows[9,8]=740

$ws3.Range("A3:I12").Value = $rows

# ------------------------------------------------------------------
# 4. Column widths on the new data sheet (best-fit approximation)
# ------------------------------------------------------------------
$ws3.Columns.Item(2).ColumnWidth = 14.42
$ws3.Columns.Item(3).ColumnWidth = 8.75
$ws3.Columns.Item(4).ColumnWidth = 14.42
$ws3.Columns.Item(5).ColumnWidth = 8.75
$ws3.Columns.Item(6).ColumnWidth = 14.42
$ws3.Columns.Item(7).ColumnWidth = 8.75
$ws3.Columns.Item(8).ColumnWidth = 14.42
$ws3.Columns.Item(9).ColumnWidth = 8.75

# ------------------------------------------------------------------
# 5. View state: "FSR_2 Input Test + Sensor 2" ends up the active/visible tab
# ------------------------------------------------------------------
$ws3.Activate()
$excel.ActiveWindow.Zoom = 71
$ws3.Range("I12").Select()
